# Remove the block of blank paragraphs (and the embedded figure) that
# sits between the "IiI_iIi... incidente." definition line and the
# trailing "Source : ..." paragraph. The definition paragraph and the
# source paragraph themselves are left untouched.

$d = $word.ActiveDocument

# Locate the paragraph that ends the "IiI_iIi est l'intensite de la
# lumiere incidente." definition - match on the stable ASCII substring
# "incidente" so accent/encoding quirks don't break the lookup.
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*IiI_iIi*incidente*") {
        $anchorPara = $p
    }
}

# Locate the following "Source : ..." paragraph (walk forward from the
# anchor so we only ever consider paragraphs after it).
$sourcePara = $anchorPara.Next()
while ($sourcePara.Range.Text -notlike "*Source*") {
    $sourcePara = $sourcePara.Next()
}

# Delete everything between the two (the blank paragraphs plus the
# paragraph holding the anchored picture), leaving both boundary
# paragraphs intact.
$deleteRange = $d.Range($anchorPara.Range.End, $sourcePara.Range.Start)
$deleteRange.Delete()
